$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (B5:AH5): round values to 2 decimal places (custom accuracy).
$ws.Range("B5").Value = 5.19
$ws.Range("C5").Value = 3.72
$ws.Range("D5").Value = 0.72
$ws.Range("E5").Value = 11.5
$ws.Range("F5").Value = 8.869999999999999
$ws.Range("G5").Value = 4.02
$ws.Range("H5").Value = 21.59
$ws.Range("I5").Value = 6.4
$ws.Range("J5").Value = 2.76
$ws.Range("K5").Value = 3.85
$ws.Range("L5").Value = 4.59
$ws.Range("M5").Value = 4.96
$ws.Range("N5").Value = 1.34
$ws.Range("O5").Value = 4.14
$ws.Range("P5").Value = 5.84
$ws.Range("Q5").Value = 3.72
$ws.Range("R5").Value = 0.67
$ws.Range("S5").Value = 0.39
$ws.Range("T5").Value = 55.67
$ws.Range("U5").Value = 11.85
$ws.Range("V5").Value = 3.82
$ws.Range("W5").Value = 7.77
$ws.Range("X5").Value = 3.97
$ws.Range("Y5").Value = 0.85
$ws.Range("Z5").Value = 9.94
$ws.Range("AA5").Value = 3.37
$ws.Range("AB5").Value = 3.11
$ws.Range("AC5").Value = 3.64
$ws.Range("AD5").Value = 4.78
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 20.15
$ws.Range("AG5").Value = 2.04
$ws.Range("AH5").Value = 4.77

# Remove row 6 entirely (data trimmed down / regenerated with 1000-row dataset
# upstream; this sheet's sample shrinks from 6 to 5 rows). This also updates
# the sheet's used-range dimension from A1:AH6 to A1:AH5 automatically.
$ws.Rows.Item(6).Delete()
